$d = $word.ActiveDocument

$replacements = @(
    @{old = "41×15=615"; new = "54×41=2214"},
    @{old = "52×62=3224"; new = "35×28=980"},
    @{old = "87×84=7308"; new = "56×25=1400"},
    @{old = "89×13=1157"; new = "72×96=6912"},
    @{old = "21×72=1512"; new = "91×19=1729"},
    @{old = "93×91=8463"; new = "21×97=2037"},
    @{old = "94×11=1034"; new = "58×43=2494"},
    @{old = "44×21=924"; new = "14×31=434"},
    @{old = "45×80=3600"; new = "74×60=4440"},
    @{old = "30×72=2160"; new = "13×70=910"},
    @{old = "26×88=2288"; new = "55×98=5390"},
    @{old = "21×74=1554"; new = "49×90=4410"},
    @{old = "52×25=1300"; new = "47×78=3666"},
    @{old = "34×14=476"; new = "20×20=400"},
    @{old = "13×75=975"; new = "60×70=4200"},
    @{old = "94×23=2162"; new = "11×45=495"},
    @{old = "95×36=3420"; new = "64×32=2048"},
    @{old = "18×25=450"; new = "74×75=5550"},
    @{old = "38×55=2090"; new = "51×71=3621"},
    @{old = "35×66=2310"; new = "49×19=931"},
    @{old = "92×71=6532"; new = "17×57=969"},
    @{old = "53×73=3869"; new = "58×90=5220"},
    @{old = "39×96=3744"; new = "86×30=2580"},
    @{old = "30×53=1590"; new = "89×75=6675"},
    @{old = "57×62=3534"; new = "45×36=1620"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2) | Out-Null
}
